# Add data for 2022-04-12:
#  - rename sheet / update "as-of" label from 04-03 to 04-04
#  - bump April (row 5) carjacking total in the 2022 column (I) from 7 to 11
#  - bump the yearly Total row (row 14) in the 2022 column (I) from 440 to 444

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet tab name: "Through 2022-04-03" -> "Through 2022-04-04"
$ws.Name = "Through 2022-04-04"

# Header label in I1: "2022 (through 04-03)" -> "2022 (through 04-04)"
$ws.Cells.Item(1, 9).Value = "2022 (through 04-04)"

# April 2022 value: 7 -> 11
$ws.Cells.Item(5, 9).Value = 11

# Total 2022 value: 440 -> 444
$ws.Cells.Item(14, 9).Value = 444
